$d = $word.ActiveDocument

function Replace-Text {
    param($range, $searchText, $replaceText)
    $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# Helper used only where the plain Find/Replace-All above would let the
# replaced run inherit formatting from an immediately preceding
# <w:hyperlink> run (an interop quirk). It finds the exact text once,
# inserts the new text after it, deletes the old text, then restores the
# run-level formatting explicitly on the freshly inserted range.
function Replace-TextPreserveFormat {
    param($range, $searchText, $replaceText, $color)
    $found = $range.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { return }
    $matchStart = $range.Start
    $matchEnd = $range.End
    $range.InsertAfter($replaceText)
    $d.Range($matchStart, $matchEnd).Delete()
    $newRange = $d.Range($matchStart, $matchStart + $replaceText.Length)
    if ($null -ne $color) { $newRange.Font.Color = $color }
}

Replace-Text $d.Content 'English' 'ඉංග්‍රීසි'
Replace-Text $d.Content 'English' 'ඉංග්‍රීසි'
Replace-Text $d.Content 'Brief' 'සංක්ෂिप्तය'
Replace-Text $d.Content 'An email sent to partners in the target country whose documents failed our verification process. It will be sent via customer.io' 'අපගේ සත්‍යාපන ක්‍රියාවලියෙන් අසමත් වූ ඔබගේ ලේඛන සඳහා ඉලක්ක ඉතා දේපල වීම් කළ සහයෝගීන්ට යැවුණු ඉමීලයක්. ඒවා customer.io හරහා යැවෙනු ඇත'
Replace-Text $d.Content 'Target audience' 'ඉලක්ක ප්‍රේක්ෂක'
Replace-Text $d.Content 'Invited partners who submitted wrong/incomplete documents' 'වැරදී හෝ නොපුරවන ලේඛන යොමුකළ ආරාධනා කළ හවුල්කරුවන්'
Replace-Text $d.Content 'Subject line' 'විෂය රේඛාව'
Replace-Text $d.Content ' — document verification failed ' ' — ලේඛන සත්‍යාපන අසමත් විය '
Replace-Text $d.Content 'Uh oh! Your documents couldn’t be verified' 'Uh oh! ඔබගේ ලේඛන සත්‍යාපනය කළ නොහැකි විය'
Replace-Text $d.Content 'Hi ' 'ආයුබෝවන් '
Replace-Text $d.Content 'We regret to inform you that your documents have failed our verification process as we found the following issues with them: ' 'අපි ඔබගේ ලේඛනවල​ පහත සඳහන් ගැටලු සොයා ගත් බැවින් ඒවා අපගේ සත්‍යාපන ක්‍රියාවලියෙන් අසමත් වී ඇති බව කනගාටුවෙන් ඔබට දන්වා සිටිමු: '
Replace-Text $d.Content 'A copy of your vaccination certificate' 'ඔබේ එන්නත් සහතිකයේ පිටපතක්'
Replace-Text $d.Content ': Document is unclear' ': ලේඛනය අපැහැදිලියි'
Replace-Text $d.Content '[Document 2]' '[ලේඛන 2]'
Replace-Text $d.Content ': [problem]' ': [ගැටලුව]'
Replace-Text $d.Content 'Please resubmit the documents above by ' 'කරුණාකර ඉහත ලේඛන '
Replace-Text $d.Content ' so we can proceed with the necessary arrangements.' ' දිනට පෙර​ නැවත ඉදිරිපත් කරන්න, එවිට අපට අවශ්‍ය විධිවිධාන සමඟ ඉදිරියට යා හැක.'
Replace-Text $d.Content 'If you have any questions, please contact your country manager, ' 'ඔබට කිසියම් ප්‍රශ්නයක් ඇත්නම්, කරුණාකර ඔබගේ රටේ කළමනාකරු, '
Replace-Text $d.Content ', at ' ', '
Replace-Text $d.Content ' (WhatsApp). ' ' (WhatsApp) හරහා අමතන්න. '

# " / Portuguese / French / Thai / Vietnamese / Spanish" -> Sinhala: this run
# immediately follows the "English" hyperlink run, so use the format-preserving helper
# (restores the original red w:color="ff0000"; the source run has no underline).
Replace-TextPreserveFormat $d.Content ' / Portuguese / French / Thai / Vietnamese / Spanish' ' / පෝරුසිය / ප්‍රංශ / තැයි / වීට්නාමීස් / ස්පාඤ්ජ නම්' 255

# Second " or " occurrence only (between [EMAIL ADDRESS] and [WHATSAPP NO]);
# the first " or " (in the "live chat ... or ... WhatsApp" sentence) must stay unchanged.
$emailRange = $d.Content
$found = $emailRange.Find.Execute('[EMAIL ADDRESS]', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $scoped = $d.Range($emailRange.End, $d.Content.End)
    Replace-Text $scoped ' or ' ' හරහා හෝ '
}
